# Apply "recolor DNA damage figures" update to the per-sample-by-cell table.
# This updates the statistics on sheet 2 (DNA_dam_score_by_sampleCelltype),
# adds a new "celltype" column (G), reorders/updates the Mural vs D2-MSN rows,
# and appends a new row for the D1/D2-Hybrid celltype.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# New header cell for column G
$ws.Range("G1").Value = "celltype"

# Data rows: term, estimate, std.error, statistic, df, p.value, celltype
$data = @(
    @{ Row = 2;  Term = "celltype4Microglia:DSM.IV.OUDOUD";      B = 0.005237570917737962;   C = 0.001765807011484599;  D = 2.96610608275617;    E = 76.34068802134549;   F = 0.004024241284508979;  G = "Microglia" },
    @{ Row = 3;  Term = "celltype4Endothelial:DSM.IV.OUDOUD";    B = 0.005117871645542164;   C = 0.001798781290339208;  D = 2.845188391178482;   E = 80.13137236548458;   F = 0.00563263916610259;   G = "Endothelial" },
    @{ Row = 4;  Term = "celltype4Oligos_Pre:DSM.IV.OUDOUD";     B = 0.004146974536677158;   C = 0.001765211947161378;  D = 2.349278534708464;   E = 76.2620154281329;    F = 0.02139995384339828;   G = "Oligos_Pre" },
    @{ Row = 5;  Term = "celltype4Interneurons:DSM.IV.OUDOUD";   B = 0.004079075694946135;   C = 0.001764168583079837;  D = 2.312180215694011;   E = 76.1132943218489;    F = 0.02347327185342079;   G = "Interneurons" },
    @{ Row = 6;  Term = "celltype4Astrocytes:DSM.IV.OUDOUD";     B = 0.002707047869883303;   C = 0.001765280720146236;  D = 1.533494270338517;   E = 76.27121195121346;   F = 0.1292919146334109;    G = "Astrocytes" },
    @{ Row = 7;  Term = "celltype4Oligos:DSM.IV.OUDOUD";         B = 0.002592599894334472;   C = 0.001764141056642705;  D = 1.469610315213902;   E = 76.08395242125513;   F = 0.1457904234221268;    G = "Oligos" },
    @{ Row = 8;  Term = "celltype4D1-MSN:DSM.IV.OUDOUD";         B = 0.001948846654841802;   C = 0.00176403634475677;   D = 1.104765590932603;   E = 76.08583330064657;   F = 0.2727425945235774;    G = "D1-MSN" },
    @{ Row = 9;  Term = "celltype4D2-MSN:DSM.IV.OUDOUD";         B = 0.0009584590446104238;  C = 0.001764036668905779;  D = 0.5433328351416584;  E = 76.0822131047389;    F = 0.5884890948944963;    G = "D2-MSN" },
    @{ Row = 10; Term = "celltype4Mural:DSM.IV.OUDOUD";          B = -0.0004488808464833061; C = 0.001836770867981745;  D = -0.2443858699569528;  E = 84.10064528662075;   F = 0.8075272301508385;    G = "Mural" },
    @{ Row = 11; Term = "celltype4D1/D2-Hybrid:DSM.IV.OUDOUD";   B = 0.0002850507650206718;  C = 0.001764248264487787;  D = 0.1615706648312506;  E = 76.12620225943856;   F = 0.8720720085822798;    G = "D1/D2-Hybrid" }
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Range("A$r").Value = $row.Term
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = $row.G
}
